# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. Insert a new "Player Info" sheet as the first sheet in the workbook.
# 2. Rewrite the MATCH_CARD_LINK columns on "ODI Batting" / "ODI Bowling"
#    as MATCH_CODE columns holding just the numeric match code instead of
#    the full scorecard URL.

$wb = $excel.ActiveWorkbook

# --- 1. New "Player Info" sheet -------------------------------------------
# Worksheets.Add() with no arguments inserts the new sheet immediately
# before the (current) first/active sheet, so this naturally lands the new
# sheet in slot 1 and pushes "ODI Batting"/"ODI Bowling" down by one.
$infoSheet = $wb.Worksheets.Add()
$infoSheet.Name = "Player Info"

$infoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $infoHeaders.Length; $i++) {
    $infoSheet.Cells.Item(1, $i + 1).Value = $infoHeaders[$i]
}

$infoRow = @("'3968", "Dean Elgar", "Left Handed", "Left Arm Orthodox")
for ($i = 0; $i -lt $infoRow.Length; $i++) {
    $infoSheet.Cells.Item(2, $i + 1).Value = $infoRow[$i]
}

# Match the bold / bordered / centered header look used by the other sheets.
$infoHeaderRange = $infoSheet.Range("A1:D1")
$infoHeaderRange.Font.Bold = $true
$infoHeaderRange.Borders.LineStyle = 1
$infoHeaderRange.HorizontalAlignment = -4108
$infoHeaderRange.VerticalAlignment = -4160

$infoSheet.Range("A1").Select()

# --- 2. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE -------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

# Rows 2..9, in order -> the numeric match code extracted from the old URL.
$battingCodes = @("3436", "3438", "3440", "3442", "3444", "3851", "4206", "4207")
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $row = $i + 2
    # Leading apostrophe forces text storage so the numeric-looking match
    # code round-trips as a string (matching MATCH_CODE's inlineStr type)
    # instead of silently becoming a numeric cell.
    $battingSheet.Cells.Item($row, 4).Value = "'" + $battingCodes[$i]
}

# --- 3. ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE --------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

# Rows 2..5, in order -> the numeric match code extracted from the old URL.
$bowlingCodes = @("3438", "3440", "3442", "3444")
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $row = $i + 2
    $bowlingSheet.Cells.Item($row, 2).Value = "'" + $bowlingCodes[$i]
}
